$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers and must be forced to
# text format first so Excel does not re-interpret/round them as numbers
# (e.g. "509.10" -> 509.1, "0.0000140" -> 1.4E-05).
$textCells = @("D5","D6","D7","D8","D10","D16","D17","D19","D20","D24","D27","D29","D31","D35","D37","D39","D40","D41","D42","D43","D45","D47")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (price, link, name, 1h volume %).
$ws.Range('D2').Value = '60.324.51'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '2.592.65'
$ws.Range('E3').Value = '  -2.79%  '
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').Value = '509.10'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').Value = '153.65'
$ws.Range('E6').Value = '  -2.20%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').Value = '  -2.63%  '
$ws.Range('D9').Value = '2.599.19'
$ws.Range('E9').Value = '  -2.33%  '
$ws.Range('D10').Value = '6.67'
$ws.Range('E10').Value = '  +5.37%  '
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('E13').Value = '  +1.70%  '
$ws.Range('D14').Value = '3.045.78'
$ws.Range('E14').Value = '  -1.78%  '
$ws.Range('D15').Value = '60.304.60'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').Value = '21.46'
$ws.Range('E16').Value = '  -1.72%  '
$ws.Range('D17').Value = '0.0000140'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('D18').Value = '2.598.27'
$ws.Range('E18').Value = '  -2.07%  '
$ws.Range('D19').Value = '4.73'
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').Value = '354.66'
$ws.Range('E20').Value = '  +1.58%  '
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('E22').Value = '  -0.76%  '
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').Value = '60.32'
$ws.Range('E24').Value = '  +0.32%  '
$ws.Range('E25').Value = '  -0.41%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('D28').Value = '0.0₃0833'
$ws.Range('E28').Value = '  -3.17%  '
$ws.Range('D29').Value = '7.33'
$ws.Range('E29').Value = '  -2.73%  '
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('D31').Value = '19.32'
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('E32').Value = '  -4.20%  '
$ws.Range('E33').Value = '  -1.07%  '
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('D35').Value = '3.98'
$ws.Range('E35').Value = '  -1.18%  '
$ws.Range('E36').Value = '  -3.06%  '
$ws.Range('D37').Value = '0.874'
$ws.Range('E38').Value = '  -2.90%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '36.08'
$ws.Range('E39').Value = '  +2.12%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '3.75'
$ws.Range('E40').Value = '  -0.57%  '
$ws.Range('D41').Value = '0.837'
$ws.Range('E41').Value = '  -2.92%  '
$ws.Range('D42').Value = '294.77'
$ws.Range('E42').Value = '  -4.40%  '
$ws.Range('D43').Value = '0.100'
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('E44').Value = '  -4.18%  '
$ws.Range('D45').Value = '0.997'
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('E46').Value = '  -4.64%  '
$ws.Range('D47').Value = '19.64'
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('E48').Value = '  -2.86%  '
$ws.Range('E49').Value = '  -1.66%  '
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('D51').Value = '1.987.43'
$ws.Range('E51').Value = '  -2.45%  '
